# Regenerate merged AHB files
#
# 1. Rename the header row: "_old" -> "_FV2404" and "_new" -> "_FV2410"
#    (the "diff" column header is left untouched).
# 2. Turn the used range A1:U92 into an Excel Table ("Table1") with an
#    AutoFilter on the header row.
# 3. Freeze the header row (split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename header cells -------------------------------------------------
$headers = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404",
    "Segment ID_FV2404", "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
    "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Convert A1:U92 into a table with an AutoFilter ----------------------
$tableRange = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $true
